{"js": "// Apply the textual edits described by the diff using the Word JavaScript\n// API (Office.js). Each change is located with Body.search() (which is\n// whitespace / nbsp tolerant) and then swapped in place with\n// Range.insertText(..., Word.InsertLocation.replace) so that only the\n// minimal changed substring is touched and the rest of each paragraph\n// (and its runs/formatting) is left alone.\n\nasync function replaceOnce(body, searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true, matchWholeWord: false }, options || {});\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1. \"_ marketing (\" -> \"-Marketing (\"\nawait replaceOnce(body, \"_ marketing (\", \"-Marketing (\");\n\n// 2. \"_Commerciaux\" -> \"-Commerciaux\"\nawait replaceOnce(body, \"_Commerciaux\", \"-Commerciaux\");\n\n// 3. \"_chef de projet\" -> \"-Chef de projet\"\nawait replaceOnce(body, \"_chef de projet\", \"-Chef de projet\");\n\n// 4. \"peuvent elles\" -> \"peuvent-elles\"\nawait replaceOnce(body, \"peuvent elles\", \"peuvent-elles\");\n\n// 5. trailing \" :\" removed after \"hors de l'entreprise ?\" (only touch the\n//    \"? :\" itself so the existing nbsp before \"?\" is left untouched)\nawait replaceOnce(body, \"? :\", \"? \");\n\n// 6. \" En fonction de l'orientation\" -> \" Oui, en fonction de l'orientation\"\nawait replaceOnce(body, \" En fonction de l\\u2019orientation\", \" Oui, en fonction de l\\u2019orientation\");\n\n// 7. \"1-2 mois\" -> \"1 \u00e0 2 mois\" (leave the rest of the sentence untouched)\nawait replaceOnce(body, \"1-2 mois\", \"1 \\u00e0 2 mois\");\n\n// 8. \"big data\" -> \"Big Data\"\nawait replaceOnce(body, \"contexte de big data pou\", \"contexte de Big Data pou\");\n\n// 9. \"peut etre retenu\" -> \"peut \u00eatre retenu\"\nawait replaceOnce(body, \"peut etre retenu\", \"peut \\u00eatre retenu\");\n", "ps1": "# Apply the textual edits described by the diff using the Word COM object\n# model. Each change is located with Range.Find/Execute (wildcard-free,\n# case-sensitive) using the minimal anchor text so only the changed\n# substring is touched and the rest of each paragraph (and the existing\n# non-breaking spaces used by the French typography in this document) is\n# left alone.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $FindText\"\n    }\n}\n\n# 1. \"_ marketing (\" -> \"-Marketing (\"\nReplace-Text \"_ marketing (\" \"-Marketing (\"\n\n# 2. \"_Commerciaux\" -> \"-Commerciaux\"\nReplace-Text \"_Commerciaux\" \"-Commerciaux\"\n\n# 3. \"_chef de projet\" -> \"-Chef de projet\"\nReplace-Text \"_chef de projet\" \"-Chef de projet\"\n\n# 4. \"peuvent elles\" -> \"peuvent-elles\"\nReplace-Text \"peuvent elles\" \"peuvent-elles\"\n\n# 5. trailing \" :\" removed after \"hors de l\u2019entreprise ?\" (only touch the\n#    \"? :\" itself so the existing nbsp before \"?\" is left untouched)\nReplace-Text \"? :\" \"? \"\n\n# 6. \" En fonction de l\u2019orientation\" -> \" Oui, en fonction de l\u2019orientation\"\nReplace-Text \" En fonction de l\u2019orientation\" \" Oui, en fonction de l\u2019orientation\"\n\n# 7. \"1-2 mois\" -> \"1 \u00e0 2 mois\" (leave the rest of the sentence untouched)\nReplace-Text \"1-2 mois\" \"1 \u00e0 2 mois\"\n\n# 8. \"big data\" -> \"Big Data\"\nReplace-Text \"contexte de big data pou\" \"contexte de Big Data pou\"\n\n# 9. \"peut etre retenu\" -> \"peut \u00eatre retenu\"\nReplace-Text \"peut etre retenu\" \"peut \u00eatre retenu\"\n"}
